$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.788.22'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '2.956.10'
$ws.Range("E3").Value = '  +2.70%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = "'199.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.30%  '

$ws.Range("D6").Value = "'595.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.75%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.89%  '

$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").Value = '2.955.68'
$ws.Range("E10").Value = '  +2.67%  '

$ws.Range("D11").Value = "'0.446"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +13.76%  '

$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").Value = '3.508.24'
$ws.Range("E13").Value = '  +3.14%  '

$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D15").Value = '76.757.20'
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("D16").Value = "'28.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.65%  '

$ws.Range("E17").Value = '  -1.28%  '

$ws.Range("D18").Value = '2.929.60'
$ws.Range("E18").Value = '  +1.89%  '

$ws.Range("D19").Value = "'13.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.53%  '

$ws.Range("D20").Value = "'8.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.93%  '

$ws.Range("D21").Value = "'371.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.38%  '

$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.13%  '

$ws.Range("D23").Value = "'2.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.25%  '

$ws.Range("D24").Value = "'72.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("D26").Value = '3.102.94'
$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("E28").Value = '  -1.06%  '

$ws.Range("D29").Value = "'0.0000107"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.05%  '

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").Value = "'8.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.83%  '

$ws.Range("E32").Value = '  -2.96%  '

$ws.Range("D33").Value = "'495.61"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.57%  '

$ws.Range("D34").Value = "'1.83"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = "'166.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.26%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = "'0.401"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +16.09%  '

$ws.Range("B38").Value = 'Cronos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D38").Value = "'0.112"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +20.50%  '

$ws.Range("D39").Value = "'20.14"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("D40").Value = "'19.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.30%  '

$ws.Range("E41").Value = '  -6.92%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").Value = "'181.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.17%  '

$ws.Range("E44").Value = '  -3.77%  '

$ws.Range("E45").Value = '  -2.32%  '

$ws.Range("D46").Value = "'40.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.54%  '

$ws.Range("E47").Value = '  -4.63%  '

$ws.Range("D48").Value = "'0.591"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.98%  '

$ws.Range("D49").Value = "'3.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.25%  '

$ws.Range("D50").Value = "'2.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.42%  '

$ws.Range("D51").Value = "'22.63"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.91%  '
